# Update the resampled-performance metrics table (Sheet1, rows 2-14) to the
# refreshed evaluation numbers (Pos pred counts + derived precision/recall/
# F1/micro/macro/balanced-accuracy figures).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E2").Value = "73/509"
$ws.Range("F2").Value = 0.8493150684931506
$ws.Range("G2").Value = 0.5904761904761905
$ws.Range("H2").Value = 0.916
$ws.Range("I2").Value = 0.697
$ws.Range("J2").Value = 0.869
$ws.Range("K2").Value = 0.806
$ws.Range("L2").Value = 0.7772642390289449

$ws.Range("E3").Value = "98/509"
$ws.Range("F3").Value = 0.9081632653061225
$ws.Range("G3").Value = 0.712
$ws.Range("H3").Value = 0.916
$ws.Range("I3").Value = 0.798
$ws.Range("J3").Value = 0.881
$ws.Range("K3").Value = 0.857
$ws.Range("L3").Value = 0.8382134387351778

$ws.Range("E4").Value = "99/509"
$ws.Range("F4").Value = 0.9090909090909091
$ws.Range("H4").Value = 0.904
$ws.Range("I4").Value = 0.773
$ws.Range("J4").Value = 0.865
$ws.Range("K4").Value = 0.838
$ws.Range("L4").Value = 0.8183790350572718

$ws.Range("E5").Value = "109/509"
$ws.Range("F5").Value = 0.8256880733944955
$ws.Range("G5").Value = 0.5660377358490566
$ws.Range("H5").Value = 0.845
$ws.Range("I5").Value = 0.672
$ws.Range("J5").Value = 0.789
$ws.Range("K5").Value = 0.758
$ws.Range("L5").Value = 0.7463393312449916

$ws.Range("E6").Value = "156/509"
$ws.Range("F6").Value = 0.8653846153846154
$ws.Range("G6").Value = 0.8333333333333334
$ws.Range("H6").Value = 0.907
$ws.Range("I6").Value = 0.849
$ws.Range("J6").Value = 0.885
$ws.Range("K6").Value = 0.878
$ws.Range("L6").Value = 0.8753280839895013

$ws.Range("E7").Value = "142/509"
$ws.Range("F7").Value = 0.9225352112676056
$ws.Range("G7").Value = 0.793939393939394
$ws.Range("H7").Value = 0.9
$ws.Range("I7").Value = 0.853
$ws.Range("J7").Value = 0.881
$ws.Range("K7").Value = 0.877
$ws.Range("L7").Value = 0.8711481007255655

$ws.Range("E8").Value = "198/509"
$ws.Range("F8").Value = 0.7777777777777778
$ws.Range("G8").Value = 0.8415300546448088
$ws.Range("H8").Value = 0.767
$ws.Range("I8").Value = 0.808
$ws.Range("J8").Value = 0.79
$ws.Range("K8").Value = 0.788
$ws.Range("L8").Value = 0.7866186858589898

$ws.Range("E9").Value = "174/509"
$ws.Range("F9").Value = 0.8908045977011494
$ws.Range("G9").Value = 0.8201058201058201
$ws.Range("H9").Value = 0.882
$ws.Range("I9").Value = 0.854
$ws.Range("J9").Value = 0.87
$ws.Range("K9").Value = 0.868
$ws.Range("L9").Value = 0.8664749284015338

$ws.Range("E10").Value = "196/509"
$ws.Range("F10").Value = 0.8061224489795918
$ws.Range("G10").Value = 0.8102564102564103
$ws.Range("H10").Value = 0.795
$ws.Range("I10").Value = 0.808
$ws.Range("J10").Value = 0.802
$ws.Range("K10").Value = 0.801
$ws.Range("L10").Value = 0.8013030685161833

$ws.Range("D11").Value = "199/509"
$ws.Range("E11").Value = "182/509"
$ws.Range("F11").Value = 0.8241758241758241
$ws.Range("G11").Value = 0.7537688442211056
$ws.Range("H11").Value = 0.802
$ws.Range("I11").Value = 0.787
$ws.Range("J11").Value = 0.795
$ws.Range("K11").Value = 0.795
$ws.Range("L11").Value = 0.7952517690493284

$ws.Range("D12").Value = "203/509"
$ws.Range("E12").Value = "188/509"
$ws.Range("F12").Value = 0.898936170212766
$ws.Range("G12").Value = 0.8325123152709359
$ws.Range("H12").Value = 0.864
$ws.Range("I12").Value = 0.864
$ws.Range("J12").Value = 0.864
$ws.Range("K12").Value = 0.864
$ws.Range("L12").Value = 0.865724242741851

$ws.Range("D13").Value = "237/509"
$ws.Range("E13").Value = "245/509"
$ws.Range("F13").Value = 0.8122448979591836
$ws.Range("G13").Value = 0.8396624472573839
$ws.Range("H13").Value = 0.738
$ws.Range("I13").Value = 0.826
$ws.Range("J13").Value = 0.791
$ws.Range("K13").Value = 0.782
$ws.Range("L13").Value = 0.7795873211896676

$ws.Range("D14").Value = "257/509"
$ws.Range("E14").Value = "311/509"
$ws.Range("F14").Value = 0.7813504823151125
$ws.Range("G14").Value = 0.9455252918287937
$ws.Range("H14").Value = 0.65
$ws.Range("I14").Value = 0.856
$ws.Range("J14").Value = 0.796
$ws.Range("K14").Value = 0.753
$ws.Range("L14").Value = 0.7366515348032858

